# Update the "想去人数" (attendance count) figures for several events.
# These numbers appear in both the "展览" sheet and the aggregate "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (rows correspond to each exhibition entry)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 108
$ws1.Range("F7").Value = 1153
$ws1.Range("F8").Value = 387
$ws1.Range("F13").Value = 399
$ws1.Range("F14").Value = 789
$ws1.Range("F16").Value = 724
$ws1.Range("F17").Value = 285

# Sheet "全部类型" (same events, different row numbers)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 108
$ws4.Range("F9").Value = 1153
$ws4.Range("F10").Value = 387
$ws4.Range("F20").Value = 399
$ws4.Range("F21").Value = 789
$ws4.Range("F23").Value = 724
$ws4.Range("F24").Value = 285
